# sn: update senegal forms
# Remove the "p_num" (order number) question from the survey sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows(7).Delete()
